$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46 (shifts the former row 46 "Vehicle Registration
# State" mapping, and everything below it, down by one row)
$ws.Rows(46).Insert()

# The insert drags column E's formatting onto the new row; the new mapping
# row only has data in C and F, so clear that stray cell completely.
$ws.Cells.Item(46, 5).Clear()

# Populate the new row with the "Vehicle Registration Non-Expiring Indicator" mapping
$ws.Range("C46").Value = "Vehicle Registration Non-Expiring Indicator"
$ws.Range("F46").Value = "/wm-req-doc:WarrantIssuedReport/j:ConveyanceRegistration[not(j:RegistrationExpirationDate)]/wm-req-ext:ConveyanceRegistrationNonExpiringIndicator"

# Match the row height used for this new mapping row
$ws.Rows(46).RowHeight = 56

# Reflect where the author's cursor ended up after typing the new XPath
[void]$ws.Range("F47").Select()
